# "Two Pointer + Hashing - Day 2"
# Mark several LeetCode problems as done (yellow highlight), reorganize the
# "Two Pointer" section by moving "Find the Duplicate Number" up next to the
# other Linked List Cycle problems, and insert a spacer row before the
# "Trapping Rain Water" / Intervals block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arrays - Hasing & Two Pointer")

$DONE_COLOR = 65535   # RGB(255,255,0) yellow - matches existing "done" style (fillId=3)

# --- Hashing section: mark first three easy problems as done ---
$ws.Range("A2").Interior.Color = $DONE_COLOR
$ws.Range("A3").Interior.Color = $DONE_COLOR
$ws.Range("A4").Interior.Color = $DONE_COLOR

# --- Two Pointer section ---
# "Linked List Cycle" (row 25) done
$ws.Range("A25").Interior.Color = $DONE_COLOR

# Clear the empty spacer cell A26 entirely (content + formatting) so the row
# disappears from the sheet rather than merely becoming blank.
$ws.Rows(26).Clear()

# "Linked List Cycle II" (row 28) done
$ws.Range("A28").Interior.Color = $DONE_COLOR

# Move "Find the Duplicate Number" up to directly follow "Linked List Cycle II":
# insert a new row at 29 and populate it with that problem's data.
$ws.Rows(29).Insert()
$ws.Range("A29").Value2 = "Find the Duplicate Number"
$ws.Range("B29").Value2 = "https://leetcode.com/problems/find-the-duplicate-number/description/"
$ws.Range("C29").Value2 = "Medium"
$ws.Range("A29").Interior.Color = $DONE_COLOR

# The original "Find the Duplicate Number" row has shifted down to row 35;
# remove it now that its content lives at row 29.
$ws.Rows(35).Delete()

# Mark the remaining Two Pointer medium problems (now rows 30-34) and
# "Next Permutation" (row 35) as done.
$ws.Range("A30").Interior.Color = $DONE_COLOR
$ws.Range("A31").Interior.Color = $DONE_COLOR
$ws.Range("A32").Interior.Color = $DONE_COLOR
$ws.Range("A33").Interior.Color = $DONE_COLOR
$ws.Range("A34").Interior.Color = $DONE_COLOR
$ws.Range("A35").Interior.Color = $DONE_COLOR

# Insert a spacer row before "Trapping Rain Water" (currently row 36),
# pushing it and the Intervals block below it down by one row. Clear the
# newly inserted row so it doesn't inherit formatting from the row above.
$ws.Rows(36).Insert()
$ws.Rows(36).Clear()

# Update the active selection to match the saved view.
$ws.Range("A9:A16").Select()
